# gsc-export-old/Coverage.xlsx : updated legacy GSC export data
#
# The "Chart" sheet holds a rolling daily window of Google Search Console
# coverage numbers (Date / Not indexed / Indexed / Impressions). The export
# rolled forward by one day: the oldest date row (2025-11-05) drops off the
# front of the window and every remaining row shifts up one place, so the
# sheet now runs 2025-11-06 .. 2026-01-31 (88 data rows) instead of
# 2025-11-05 .. 2026-01-31 (89 data rows).
# i.e. simply delete the sheet's row 2 (the first data row under the header)
# and let Excel shift everything below it upward. The header row (row 1) and
# all the other worksheets (Critical issues / Non-critical issues / Metadata)
# keep their own data untouched - Excel's shared-string table compaction
# takes care of renumbering their <c t="s"> references automatically.

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

# Delete the whole row 2 (2025-11-05), shifting rows 3:89 up to 2:88.
$chart.Rows(2).Delete()
